# Update "想去人数" (want-to-go count) values that changed between scrapes.
# Sheet "展览" (Exhibition)
$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 202
$wsExhibition.Range("F3").Value = 550
$wsExhibition.Range("F4").Value = 48
$wsExhibition.Range("F7").Value = 36
$wsExhibition.Range("F8").Value = 30
$wsExhibition.Range("F9").Value = 420
$wsExhibition.Range("F10").Value = 3487
$wsExhibition.Range("F11").Value = 51

# Sheet "演出" (Performance)
$wsPerformance = $wb.Worksheets.Item("演出")
$wsPerformance.Range("F2").Value = 96
$wsPerformance.Range("F3").Value = 46

# Sheet "全部类型" (All Types) - consolidated view of the two sheets above
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 202
$wsAll.Range("F3").Value = 96
$wsAll.Range("F4").Value = 550
$wsAll.Range("F5").Value = 48
$wsAll.Range("F8").Value = 36
$wsAll.Range("F9").Value = 30
$wsAll.Range("F10").Value = 420
$wsAll.Range("F11").Value = 3487
$wsAll.Range("F12").Value = 51
$wsAll.Range("F13").Value = 46
